# Append a new tracker row (row 2) to the active sheet, matching the
# "Update tracker data 2025-07-22 17:46:25" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "G1"
$ws.Range("B2").Value = "Test1"
$ws.Range("C2").Value = "Daily"
$ws.Range("D2").Value = 1

# DateAdded: 2025-07-22 stored as the Excel date serial (45860), no time
# component. Apply the lowercase format first (registers numFmtId 164)
# then switch to the uppercase variant actually used by the cell style
# (numFmtId 165), matching both numFmt entries captured in the diff.
$ws.Range("E2").Value = 45860
$ws.Range("E2").NumberFormat = "yyyy-mm-dd"
$ws.Range("E2").NumberFormat = "YYYY-MM-DD"

$ws.Range("F2").Value = 30
